# Update market-price / profit columns (H:N) across several worksheets
# (ALC, ARM, BSM, CUL, LTW, WVR) to reflect refreshed Kraken market data,
# as produced by the scheduled Sheets runner.
$wb = $excel.ActiveWorkbook

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3297.8333
$ws.Range("I43").Value = 3330.3333
$ws.Range("J43").Value = 3265.3333
$ws.Range("K43").Value = 3330.3333
$ws.Range("L43").Value = 3265.3333
$ws.Range("M43").Value = -3261.3333
$ws.Range("N43").Value = -3403.3333

# ALC row 75
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 20000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 20000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 20000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -21872

# ALC row 78
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 20000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 20000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 60000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -69360

# ALC row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 239.75
$ws.Range("I99").Value = 259.33334
$ws.Range("J99").Value = 181
$ws.Range("K99").Value = 778.0000200000001
$ws.Range("L99").Value = 543
$ws.Range("M99").Value = 719.9999799999999
$ws.Range("N99").Value = -3539

# ALC row 105
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 21000
$ws.Range("I105").Value = 10000
$ws.Range("J105").Value = 26500
$ws.Range("K105").Value = 10000
$ws.Range("L105").Value = 26500
$ws.Range("M105").Value = -6506
$ws.Range("N105").Value = -33488

# ALC row 109
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 69997.5
$ws.Range("J109").Value = 69997.5
$ws.Range("L109").Value = 69997.5
$ws.Range("N109").Value = -72771.5

# ALC row 114
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 99995
$ws.Range("J114").Value = 99995
$ws.Range("L114").Value = 99995
$ws.Range("N114").Value = -108673

# ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 194.75
$ws.Range("I118").Value = 194.75
$ws.Range("K118").Value = 584.25
$ws.Range("M118").Value = 1072.75

# ALC row 124
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H124").Value = 66497.5
$ws.Range("J124").Value = 66497.5
$ws.Range("L124").Value = 66497.5
$ws.Range("N124").Value = -76317.5

# ALC row 126
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 66497.5
$ws.Range("J126").Value = 66497.5
$ws.Range("L126").Value = 66497.5
$ws.Range("N126").Value = -76377.5

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2949.5
$ws.Range("I129").Value = 5000
$ws.Range("J129").Value = 899
$ws.Range("K129").Value = 15000
$ws.Range("L129").Value = 2697
$ws.Range("M129").Value = -10000
$ws.Range("N129").Value = -12697

# ALC row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 94976.5
$ws.Range("J130").Value = 94976.5
$ws.Range("L130").Value = 94976.5
$ws.Range("N130").Value = -105016.5

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3334.524
$ws.Range("I132").Value = 3334.524
$ws.Range("K132").Value = 10003.572
$ws.Range("M132").Value = -7473.572

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3853.0908
$ws.Range("I138").Value = 2450
$ws.Range("J138").Value = 3993.4
$ws.Range("K138").Value = 7350
$ws.Range("L138").Value = 11980.2
$ws.Range("M138").Value = -2210
$ws.Range("N138").Value = -22260.2

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1391.6578
$ws.Range("I32").Value = 1357.8889
$ws.Range("K32").Value = 1357.8889
$ws.Range("M32").Value = -1070.8889

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 864.5294
$ws.Range("I63").Value = 893.5625
$ws.Range("J63").Value = 400
$ws.Range("K63").Value = 893.5625
$ws.Range("L63").Value = 400
$ws.Range("M63").Value = -207.5625
$ws.Range("N63").Value = -1772

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 864.5294
$ws.Range("I66").Value = 893.5625
$ws.Range("J66").Value = 400
$ws.Range("K66").Value = 4467.8125
$ws.Range("L66").Value = 2000
$ws.Range("M66").Value = -1035.8125
$ws.Range("N66").Value = -8864

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1624.375
$ws.Range("I88").Value = 1416.6666
$ws.Range("J88").Value = 2247.5
$ws.Range("K88").Value = 1416.6666
$ws.Range("L88").Value = 2247.5
$ws.Range("M88").Value = -1010.6666
$ws.Range("N88").Value = -3059.5

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1624.375
$ws.Range("I91").Value = 1416.6666
$ws.Range("J91").Value = 2247.5
$ws.Range("K91").Value = 1416.6666
$ws.Range("L91").Value = 2247.5
$ws.Range("M91").Value = -12.66660000000002
$ws.Range("N91").Value = -5055.5

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1585.5714
$ws.Range("J102").Value = 1499.5
$ws.Range("L102").Value = 1499.5
$ws.Range("N102").Value = -4743.5

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 77000
$ws.Range("J133").Value = 77000
$ws.Range("L133").Value = 77000
$ws.Range("N133").Value = -82060

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2521.087
$ws.Range("I94").Value = 2170.75
$ws.Range("K94").Value = 2170.75
$ws.Range("M94").Value = -1719.75

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3335
$ws.Range("I105").Value = 3043.75
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 3043.75
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -1296.75
$ws.Range("N105").Value = -7994

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 729.875
$ws.Range("I12").Value = 349.8
$ws.Range("J12").Value = 1363.3334
$ws.Range("K12").Value = 1049.4
$ws.Range("L12").Value = 4090.0002
$ws.Range("M12").Value = -876.4000000000001
$ws.Range("N12").Value = -4436.0002

# CUL row 82
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 2999
$ws.Range("I82").Value = 2999
$ws.Range("K82").Value = 8997
$ws.Range("M82").Value = -8591

# CUL row 85
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 2999
$ws.Range("I85").Value = 2999
$ws.Range("K85").Value = 8997
$ws.Range("M85").Value = -7593

# CUL row 124
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 3000
$ws.Range("I124").Value = 3000
$ws.Range("K124").Value = 9000
$ws.Range("M124").Value = -4090

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1327.6666
$ws.Range("I139").Value = 1589.8
$ws.Range("J139").Value = 1000
$ws.Range("K139").Value = 4769.4
$ws.Range("L139").Value = 3000
$ws.Range("M139").Value = 370.6000000000004
$ws.Range("N139").Value = -13280

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1095
$ws.Range("I22").Value = 1018.8333
$ws.Range("K22").Value = 1018.8333
$ws.Range("M22").Value = -723.8333

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1095
$ws.Range("I27").Value = 1018.8333
$ws.Range("K27").Value = 1018.8333
$ws.Range("M27").Value = -911.8333

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 25100.375
$ws.Range("I81").Value = 28543.285
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 57086.57
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -56025.57
$ws.Range("N81").Value = -4122

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 25100.375
$ws.Range("I84").Value = 28543.285
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 285432.85
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -280128.85
$ws.Range("N84").Value = -20608

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1826.9
$ws.Range("I122").Value = 2035.625
$ws.Range("K122").Value = 6106.875
$ws.Range("M122").Value = -3656.875
